# Apply the "cryptos" price/volume refresh described by the commit:
# "Updated cryptos list on Tue May  7 09:53:45 UTC 2024 with GitHub Actions"
#
# Column D ("Price") holds text values that often look numeric (e.g. "594.79",
# "0.0367"), so they are written with a leading apostrophe to force Excel to
# store them as text instead of silently converting them to numbers (which
# would also mangle values such as the PEPE subscript price).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "'64.221.11"
$ws.Range("E2").Value = "  -1.69%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "'3.117.76"
$ws.Range("E3").Value = "  -2.72%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.02%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "'594.79"
$ws.Range("E5").Value = "  -0.54%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "'158.25"
$ws.Range("E6").Value = "  +2.64%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.02%  "

# --- Row 8: XRP ---
$ws.Range("E8").Value = "  -0.06%  "

# --- Row 9: LidoStakedEther ---
$ws.Range("D9").Value = "'3.116.47"
$ws.Range("E9").Value = "  -2.69%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -4.95%  "

# --- Row 11: Toncoin ---
$ws.Range("E11").Value = "  -3.00%  "

# --- Row 12: Cardano ---
$ws.Range("E12").Value = "  -3.92%  "

# --- Row 13: Avalanche ---
$ws.Range("D13").Value = "'37.26"
$ws.Range("E13").Value = "  -5.61%  "

# --- Row 14: ShibaInu ---
$ws.Range("E14").Value = "  -5.80%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
$ws.Range("D15").Value = "'3.632.82"
$ws.Range("E15").Value = "  -2.73%  "

# --- Row 16: TRON ---
$ws.Range("E16").Value = "  -1.48%  "

# --- Row 17: Polkadot ---
$ws.Range("E17").Value = "  -2.42%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "'64.136.51"
$ws.Range("E18").Value = "  -1.40%  "

# --- Row 19: WrappedEther ---
$ws.Range("D19").Value = "'3.119.54"
$ws.Range("E19").Value = "  -2.58%  "

# --- Row 20: BitcoinCash ---
$ws.Range("D20").Value = "'478.21"
$ws.Range("E20").Value = "  -1.36%  "

# --- Row 21: Chainlink ---
$ws.Range("D21").Value = "'14.52"
$ws.Range("E21").Value = "  -4.18%  "

# --- Row 22: Polygon ---
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  -7.51%  "

# --- Row 23: Uniswap ---
$ws.Range("D23").Value = "'7.58"
$ws.Range("E23").Value = "  -4.49%  "

# --- Row 24: Fetch.AI ---
$ws.Range("E24").Value = "  +1.81%  "

# --- Row 25: InternetComputer(DFINITY) ---
$ws.Range("D25").Value = "'13.01"
$ws.Range("E25").Value = "  -6.85%  "

# --- Row 26: Litecoin ---
$ws.Range("D26").Value = "'81.45"
$ws.Range("E26").Value = "  -2.66%  "

# --- Row 27: RenderToken ---
$ws.Range("D27").Value = "'10.62"
$ws.Range("E27").Value = "  +7.34%  "

# --- Row 28: Dai ---
$ws.Range("E28").Value = "  -0.29%  "

# --- Row 29: NEARProtocol ---
$ws.Range("D29").Value = "'7.62"
$ws.Range("E29").Value = "  +1.34%  "

# --- Row 30: PancakeSwap ---
$ws.Range("E30").Value = "  -2.69%  "

# --- Row 31: FirstDigitalUSD ---
$ws.Range("E31").Value = "  -0.17%  "

# --- Row 32: ImmutableX ---
$ws.Range("E32").Value = "  -3.77%  "

# --- Row 33: Hedera ---
$ws.Range("E33").Value = "  -6.29%  "

# --- Row 34: EthereumClassic ---
$ws.Range("D34").Value = "'27.36"
$ws.Range("E34").Value = "  -4.28%  "

# --- Row 35: PEPE ---
$ws.Range("D35").Value = "'0.0" + [char]0x2083 + "0849"
$ws.Range("E35").Value = "  -5.89%  "

# --- Row 36: Mantle ---
$ws.Range("E36").Value = "  -2.27%  "

# --- Row 37: now dogwifhat (was Filecoin) ---
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.32"
$ws.Range("E37").Value = "  -7.63%  "

# --- Row 38: now Filecoin (was dogwifhat) ---
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'6.05"
$ws.Range("E38").Value = "  -5.05%  "

# --- Row 39: Stacks ---
$ws.Range("D39").Value = "'2.26"
$ws.Range("E39").Value = "  -4.99%  "

# --- Row 40: OKB ---
$ws.Range("D40").Value = "'51.02"
$ws.Range("E40").Value = "  -0.86%  "

# --- Row 41: now Cosmos (was Bittensor) ---
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'9.18"
$ws.Range("E41").Value = "  -3.16%  "

# --- Row 42: now Bittensor (was Cosmos) ---
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'448.09"
$ws.Range("E42").Value = "  -6.55%  "

# --- Row 43: TheGraph ---
$ws.Range("E43").Value = "  -2.81%  "

# --- Row 44: now Arweave (was VeChain) ---
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "'40.88"
$ws.Range("E44").Value = "  +5.73%  "

# --- Row 45: now VeChain (was Arweave) ---
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0367"
$ws.Range("E45").Value = "  -4.46%  "

# --- Row 46: Kaspa ---
$ws.Range("E46").Value = "  +0.34%  "

# --- Row 47: Maker ---
$ws.Range("D47").Value = "'2.834.03"
$ws.Range("E47").Value = "  -4.45%  "

# --- Row 48: Monero ---
$ws.Range("D48").Value = "'130.74"
$ws.Range("E48").Value = "  -0.91%  "

# --- Row 49: InjectiveProtocol ---
$ws.Range("D49").Value = "'25.99"
$ws.Range("E49").Value = "  +0.76%  "

# --- Row 51: ThetaToken ---
$ws.Range("E51").Value = "  -3.18%  "
